$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38: PMAfactor
$ws.Range("A38").Value = "PMAfactor"
$ws.Range("B38").Value = 1
$ws.Range("C38").Value = 1
$ws.Range("D38").Value = 0
$ws.Range("E38").Value = 1
$ws.Range("F38").Value = "no"
$ws.Range("G38").Value = "PMA_factor"

# Row 39: kLacTetRTup1
$ws.Range("A39").Value = "kLacTetRTup1"
$ws.Range("B39").Value = 0.01
$ws.Range("C39").Value = 10
$ws.Range("D39").Value = 1
$ws.Range("E39").Value = 98
$ws.Range("F39").Value = "yes"
$ws.Range("G39").Value = "k_{LacTetRTup1}"

# Row 40: mufactor
$ws.Range("A40").Value = "mufactor"
$ws.Range("B40").Value = 1
$ws.Range("C40").Value = 1
$ws.Range("D40").Value = 0
$ws.Range("E40").Value = 1
$ws.Range("F40").Value = "no"
$ws.Range("G40").Value = "mu_factor"

# Apply the same style used by surrounding rows (s="2") to the new rows
$ws.Range("A38:G40").Style = $ws.Range("A36").Style

# Update the view: scroll so row 25 is the top-left visible row, and move
# the active selection to A41 (just below the newly-added data).
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Range("A41").Select()

$wb.Windows.Item(1).WindowState = $wb.Windows.Item(1).WindowState
